$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: I1 = "Other found locations"
$ws.Range("I1").Value = "Other found locations"

# Row 2: F2/G2 become "not found" / "N/A"; new I2 = "" (empty string)
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"
$ws.Range("I2").Value = ""

# Row 3: E3 authors list gets extra spacing (triple space separators); new I3 = "_PMC"
$ws.Range("E3").Value = '[Hussein N.%Ali%NULL%1,   Sherko S.%Niranji%sherko.subhan@garmian.edu.krd%1,   Sirwan M. A.%Al‐Jaf%NULL%2,   Sirwan M. A.%Al‐Jaf%NULL%0]'
$ws.Range("I3").Value = "_PMC"

# Row 4: E4 authors list gets extra spacing (triple space separators); new I4 = "_PMC"
$ws.Range("E4").Value = '[Ana Valesca Fernandes Gilson%Silva%NULL%1,   Diego%Menezes%NULL%1,   Filipe Romero Rebello%Moreira%NULL%1,   Octávio Alcântara%Torres%NULL%1,   Paula Luize Camargos%Fonseca%NULL%1,   Rennan Garcias%Moreira%NULL%1,   Hugo José%Alves%NULL%1,   Vivian Ribeiro%Alves%NULL%1,   Tânia Maria de Resende%Amaral%NULL%1,   Adriano Neves%Coelho%NULL%1,   Júlia Maria%Saraiva Duarte%NULL%1,   Augusto Viana%da Rocha%NULL%1,   Luiz Gonzaga Paula%de Almeida%NULL%1,   João Locke Ferreira%de Araújo%NULL%1,   Hilton Soares%de Oliveira%NULL%1,   Nova Jersey Cláudio%de Oliveira%NULL%1,   Camila%Zolini%NULL%1,   Jôsy Hubner%de Sousa%NULL%1,   Elizângela Gonçalves%de Souza%NULL%1,   Rafael Marques%de Souza%NULL%1,   Luciana de Lima%Ferreira%NULL%1,   Alexandra%Lehmkuhl Gerber%NULL%1,   Ana Paula de Campos%Guimarães%NULL%1,   Paulo Henrique Silva%Maia%NULL%1,   Fernanda Martins%Marim%NULL%1,   Lucyene%Miguita%NULL%1,   Cristiane Campos%Monteiro%NULL%1,   Tuffi Saliba%Neto%NULL%1,   Fabrícia Soares Freire%Pugêdo%NULL%1,   Daniel Costa%Queiroz%NULL%1,   Damares Nigia Alborguetti Cuzzuol%Queiroz%NULL%1,   Luciana Cunha%Resende-Moreira%NULL%1,   Franciele Martins%Santos%NULL%1,   Erika Fernanda Carlos%Souza%NULL%1,   Carolina Moreira%Voloch%NULL%1,   Ana Tereza%Vasconcelos%NULL%1,   Renato Santana%de Aguiar%NULL%1,   Renan Pedra%de Souza%NULL%1]'
$ws.Range("I4").Value = "_PMC"
